$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing "productID..unitId"
# header/data columns one position to the right (B..M), and set the new first
# column's header to "id" (matching the product-import-template update that
# adds an "id" column ahead of "productID").
$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "id"

# The newly inserted column should carry the same width that column A had
# before the insert (~19.1 "characters", matching the column it displaced to
# B). ColumnWidth is expressed/rounded in pixel-quantized character units, so
# 19.0 is the closest input that reproduces that width.
$ws.Columns.Item(1).ColumnWidth = 19

# Move/park the active selection on the new id column's first data row, as in
# the updated template.
$ws.Range("A2").Select() | Out-Null
